$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.982.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.510.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.501.45'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.93%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.88%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.187'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.95%  '
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.10'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.073.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.507.77'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.937.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '542.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +14.17%  '
$ws.Range('E22').Value = '  -2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '93.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.88'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('E34').Value = '  -3.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '562.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '37.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.61%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0765'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('E42').Value = '  -3.47%  '
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('E44').Value = '  +5.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.222.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '138.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.67%  '
